# Apply the "new error log for windows error" edits to the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C: user_name, changed on every data row (2-16) ---
$ws.Range("C2").Value = "Maaya Inoue"
$ws.Range("C3").Value = "Maaya Inoue"
$ws.Range("C4").Value = "Maaya Inoue"
$ws.Range("C5").Value = "Maaya Inoue"
$ws.Range("C6").Value = "Maaya Inoue"
$ws.Range("C7").Value = "Maaya Inoue"
$ws.Range("C8").Value = "Maaya Inoue"
$ws.Range("C9").Value = "Maaya Inoue"
$ws.Range("C10").Value = "Maaya Inoue"
$ws.Range("C11").Value = "Maaya Inoue"
$ws.Range("C12").Value = "Maaya Inoue"
$ws.Range("C13").Value = "Maaya Inoue"
$ws.Range("C14").Value = "Maaya Inoue"
$ws.Range("C15").Value = "Maaya Inoue"
$ws.Range("C16").Value = "Maaya Inoue"

# --- Column J: capimg screenshot filenames ---
$ws.Range("J2").Value  = "bdot20240415_141954/1.png"
$ws.Range("J3").Value  = "bdot20240415_141954/2.png"
$ws.Range("J4").Value  = "bdot20240415_141954/3.png"
$ws.Range("J5").Value  = "bdot20240415_141954/4.png"
$ws.Range("J6").Value  = "bdot20240415_141954/5.png"
$ws.Range("J7").Value  = "bdot20240415_141954/5.png"
$ws.Range("J8").Value  = "bdot20240415_141954/6.png"
$ws.Range("J9").Value  = "bdot20240415_141954/7.png"
$ws.Range("J10").Value = "bdot20240415_141954/8.png"
$ws.Range("J11").Value = "bdot20240415_141954/9.png"
$ws.Range("J12").Value = "bdot20240415_141954/10.png"
$ws.Range("J13").Value = "bdot20240415_141954/1.png"
$ws.Range("J14").Value = "bdot20240415_141954/2.png"
$ws.Range("J15").Value = "bdot20240415_141954/3.png"
$ws.Range("J16").Value = "bdot20240415_141954/11.png"

# --- Column K: explanation text (rewritten/reshuffled steps) ---
$ws.Range("K2").Value  = "「スタート」ボタンをクリックする"
$ws.Range("K3").Value  = "メニューから「設定」アイコンをクリックする"
$ws.Range("K4").Value  = "左側のメニューからWindows Updateをクリックし、Windows Update画面に移動する"
$ws.Range("K5").Value  = "0x80240fff エラー"
$ws.Range("K6").Value  = "デスクトップ画面の左下にある「スタート」ボタンを右クリックする"
$ws.Range("K7").Value  = "メニューからターミナル(管理者)をクリックする"
$ws.Range("K8").Value  = "ユーザーアカウント制御と表示されているウィンドウが開いたことを確認する"
$ws.Range("K9").Value  = "PowerShellウィンドウに start-transcript と入力し、[Enter]キーを押す"
$ws.Range("K10").Value = "wuauclt.exe /resetauthorization /detectnow と入力し、[Enter]キーを押す"
$ws.Range("K11").Value = "netsh winhttp show proxy と入力し、[Enter]キーを押す"
$ws.Range("K12").Value = "netsh winhttp reset proxy と入力し、[Enter]キーを押す"
$ws.Range("K13").Value = "「スタート」ボタンをクリックする"
$ws.Range("K14").Value = "メニューから「設定」アイコンをクリックする"
$ws.Range("K15").Value = "左側のメニューからWindows Updateをクリックし、Windows Update画面に移動する"
$ws.Range("K16").Value = "「更新プログラムのチェック」ボタンをクリックする"

# --- Column B: type ---
# Row 5 becomes the new error row, row 7 reverts back to a normal operation row.
$ws.Range("B5").Value = "error"
$ws.Range("B7").Value = "operation"

# --- Columns L/M: error_type / error_content ---
# The error row moves from row 7 to row 5, so populate row 5 and clear row 7.
$ws.Range("L5").Value = "Error W"
$ws.Range("M5").Value = " エラーの Windows"
$ws.Range("L7").Value = ""
$ws.Range("M7").Value = ""
